# "Mise à jour de l'application" — add a new attendance column (BC) for the
# next training-session date (2025-09-30, Excel serial 45930), filling in
# each player's status for that date, mirroring the formatting already used
# by the previous date column (BB).
#
# Row 12 (Karim Belmahi) has no data past column AX in the source sheet (the
# player's tracked history ends earlier), so it intentionally receives no
# new BC cell — matching the source edit exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New session date header.
$ws.Range("BC1").Value = 45930
$ws.Range("BB1").Copy()
$ws.Range("BC1").PasteSpecial(-4122)   # xlPasteFormats — copy BB1's style (date, centered) onto BC1

# Per-player attendance status for the new date column.
$statuses = [ordered]@{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "RH"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "B"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $statuses.Keys) {
    $value = $statuses[$row]
    # Set the value first so the sheet's formulas (COUNTA/COUNTIF over
    # K:VQ-style ranges) recalc against the final content of the cell.
    $ws.Range("BC$row").Value = $value
    # Then copy over BB's formatting (same style used for every other
    # attendance cell in that row) without touching the value again.
    $ws.Range("BB$row").Copy()
    $ws.Range("BC$row").PasteSpecial(-4122)
}

# Restore the active selection to match the post-edit workbook state.
$ws.Range("BF8").Select() | Out-Null
